# Commit: "ajout des infos dans les contrats 25 aout 2025"
#
# Replace the representative's name "Monsieur El Hadj Mamadou FAYE" with
# "Madame Jenny MVOU" in the "ENTRE LES SOUSSIGNES" clause, and make the
# single space that follows the name its own bold run (matching the
# target OOXML run layout), while leaving the following "l'" / "effet des
# présentes, " runs untouched.

$d = $word.ActiveDocument

# --- Step 1: swap the name, keeping the existing (bold) run formatting ---
$oldName = "Monsieur El Hadj Mamadou FAYE"
$newName = "Madame Jenny MVOU"

$text = $d.Content.Text
$idx = $text.IndexOf($oldName)
if ($idx -lt 0) {
    throw "Could not find '$oldName' in the document"
}
$nameRange = $d.Range($idx, $idx + $oldName.Length)
$nameRange.Text = $newName

# --- Step 2: make the space right after the new name its own bold run ---
$text = $d.Content.Text
$afterName = $text.IndexOf($newName) + $newName.Length
$spaceRange = $d.Range($afterName, $afterName + 1)
$spaceRange.Font.Bold = 1

# --- Step 3: restore the original run boundary between "ayant pleins
#     pouvoirs à " and the following "l'" run. Editing step 1 above causes
#     the engine to coalesce every subsequent identically-formatted run in
#     the paragraph ("ayant pleins pouvoirs à ", "l'" and
#     "effet des présentes, ") into a single run; nudging the formatting of
#     "l'" (set then restore Bold) forces it back into its own run without
#     altering its visible text or formatting. ---
$text = $d.Content.Text
$marker = "l" + [char]0x2019
$lIdx = $text.IndexOf($marker, $afterName)
if ($lIdx -ge 0) {
    $lRange = $d.Range($lIdx, $lIdx + $marker.Length)
    $lRange.Font.Bold = 1
    $lRange.Font.Bold = 0
}
